# "final commit to gojek branch"
#
# Swaps the Register/SignIn sample test data from the "Gopi Sharma"
# persona to the "Gaurav Arora" persona (new e-mails, new surname-based
# password-confirmation strings, new mobile numbers) and clears the
# stale Results column ("PASS"/"SKIP") on the Register and SignIn sheets
# so the suite looks freshly un-run. Also nudges each sheet's saved
# cursor/selection, matching what a human would leave behind after
# making these edits directly in Excel.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Register sheet (column-by-column, top-to-bottom -- same order a person
# editing the grid in Excel would type the replacement values in)
# ---------------------------------------------------------------------
$register = $wb.Worksheets.Item("Register")

$register.Range("A2").Value = "arora.gaurav74@gmail.com"
$register.Range("A3").Value = "arora.gaurav75@gmail.com"

$register.Range("E2").Value = "Gaurav"
$register.Range("E3").Value = "Gaurav"

$register.Range("F2").Value = "Arora74"
$register.Range("F3").Value = "Arora75"

$register.Range("G2").Value = 1234567898
$register.Range("G3").Value = 1234567899

$register.Range("J2").Value = ""
$register.Range("J3").Value = ""

# ---------------------------------------------------------------------
# SignIn sheet
# ---------------------------------------------------------------------
$signIn = $wb.Worksheets.Item("SignIn")

$signIn.Range("A2").Value = "arora.gaurav74@gmail.com"
$signIn.Range("A3").Value = "arora.gaurav75@gmail.com"

$signIn.Range("F2").Value = ""
$signIn.Range("F3").Value = ""

# ---------------------------------------------------------------------
# Selections left behind on each sheet (Test Cases activated last so it
# stays the book's active tab, matching the saved file).
# ---------------------------------------------------------------------
$signIn.Activate() | Out-Null
$signIn.Range("F2:F3").Select() | Out-Null

$register.Activate() | Out-Null
$register.Range("J2:J3").Select() | Out-Null

$testCases = $wb.Worksheets.Item("Test Cases")
$testCases.Activate() | Out-Null
$testCases.Range("D3:D4").Select() | Out-Null

$signIn.Activate() | Out-Null
